# Update "想去人数" (number of people interested) values for a few events
# on the "展览" and "全部类型" sheets, reflecting a newer data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 802
$ws1.Range("F4").Value = 274
$ws1.Range("F5").Value = 921
$ws1.Range("F6").Value = 2227

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 802
$ws4.Range("F4").Value = 274
$ws4.Range("F7").Value = 921
$ws4.Range("F8").Value = 2227
